$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 27, shifting existing rows 27-71 down to 28-72.
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27. Columns A,B,C,E,F,G,H,I,N,O,Q,R keep the
# same values as the row that used to be row 27 (now row 28); columns
# D,J,K,L,M,P take on new values for this record.
$ws.Cells.Item(27, 1).Value = 4
$ws.Cells.Item(27, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(27, 3).Value = "Los Lagos"
$ws.Cells.Item(27, 4).Value = 44665
$ws.Cells.Item(27, 5).Value = 10
$ws.Cells.Item(27, 6).Value = 100112031
$ws.Cells.Item(27, 7).Value = "Poroto verde"
$ws.Cells.Item(27, 8).Value = "Magnum"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 35
$ws.Cells.Item(27, 11).Value = 28000
$ws.Cells.Item(27, 12).Value = 28000
$ws.Cells.Item(27, 13).Value = 28000
$ws.Cells.Item(27, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(27, 15).Value = "Región Metropolitana"
$ws.Cells.Item(27, 16).Value = 1120
$ws.Cells.Item(27, 17).Value = 25
$ws.Cells.Item(27, 18).Value = "Hortaliza"
